$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.149.48"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "2.245.82"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.02"
$ws.Range("E5").Value = "  -1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.89"
$ws.Range("E6").Value = "  -2.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  +1.19%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.95"
$ws.Range("E10").Value = "  -3.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0820"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.387.85"
$ws.Range("E14").Value = "  +6.60%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.587.00"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.834"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.65"
$ws.Range("E17").Value = "  -2.89%  "

$ws.Range("D18").Value = "44.062.40"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "0.0₃0977"
$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.22"
$ws.Range("E20").Value = "  -5.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.39"
$ws.Range("E21").Value = "  +1.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.46"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.42"
$ws.Range("E23").Value = "  +2.13%  "

$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").Value = "  -1.38%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  +1.51%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.90"
$ws.Range("E29").Value = "  +3.15%  "

$ws.Range("E30").Value = "  +2.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.09"
$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.36"
$ws.Range("E32").Value = "  -3.33%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0806"
$ws.Range("E33").Value = "  -2.55%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.34"
$ws.Range("E34").Value = "  +4.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.61"
$ws.Range("E35").Value = "  -3.06%  "

$ws.Range("E36").Value = "  +2.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.109"
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("E38").Value = "  -6.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.73"
$ws.Range("E39").Value = "  -5.55%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.88"
$ws.Range("E40").Value = "  -4.14%  "

$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.40"
$ws.Range("E41").Value = "  -5.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0299"
$ws.Range("E42").Value = "  -2.20%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "1.730.00"
$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "83.66"
$ws.Range("E45").Value = "  +4.69%  "

$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.39"
$ws.Range("E47").Value = "  -0.88%  "

$ws.Range("E48").Value = "  -4.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.15"
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.83"
$ws.Range("E50").Value = "  -2.31%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "68.61"
$ws.Range("E51").Value = "  -5.50%  "
